$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.154.95'
$ws.Range('E2').Value = '  +0.30%  '
$ws.Range('D3').Value = '3.721.52'
$ws.Range('E3').Value = '  -1.26%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '0.999'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '598.44'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '168.15'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.46%  '
$ws.Range('D7').Value = '3.717.98'
$ws.Range('E7').Value = '  -1.26%  '
$ws.Range('E8').Value = '  -0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.535'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.35%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.167'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.67%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.25'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.70%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.461'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -1.12%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '38.36'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +0.82%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.0000245'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +1.00%  '
$ws.Range('D15').Value = '4.341.97'
$ws.Range('E15').Value = '  -1.24%  '
$ws.Range('D16').Value = '3.718.64'
$ws.Range('E16').Value = '  -1.21%  '
$ws.Range('D17').Value = '68.027.37'
$ws.Range('E17').Value = '  +0.07%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '7.36'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +2.15%  '
$ws.Range('E19').Value = '  -0.65%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '17.35'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +9.19%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '490.03'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.96%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.28'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -0.24%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.728'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '84.87'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  -0.68%  '
$ws.Range('B25').Value = 'Fetch.AI'
$ws.Range('C25').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.33'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.04%  '
$ws.Range('B26').Value = 'PEPE'
$ws.Range('C26').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.0000142'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +3.57%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '12.36'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.69%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '10.15'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +0.61%  '
$ws.Range('E29').Value = '  +0.12%  '
$ws.Range('E30').Value = '  +0.68%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '7.86'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +4.46%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '2.38'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -1.49%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '31.54'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  -2.98%  '
$ws.Range('D34').Value = '3.860.39'
$ws.Range('E34').Value = '  -1.26%  '
$ws.Range('E35').Value = '  -0.74%  '
$ws.Range('D36').Value = '3.664.43'
$ws.Range('E36').Value = '  -1.20%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.997'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  -0.11%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '1.01'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '5.87'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +1.57%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.133'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -2.24%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.325'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +0.79%  '
$ws.Range('B42').Value = 'Bittensor'
$ws.Range('C42').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '431.19'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  -2.70%  '
$ws.Range('B43').Value = 'OKB'
$ws.Range('C43').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '48.81'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  -0.31%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.97'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  -0.81%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '2.86'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  -0.28%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '8.46'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +2.06%  '
$ws.Range('B47').Value = 'Arweave'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7XWg41D1+arweave-ar'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '40.84'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.69%  '
$ws.Range('B48').Value = 'USDe'
$ws.Range('C48').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.00'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +0.03%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '141.07'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +2.58%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.0353'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.13%  '
$ws.Range('D51').Value = '2.758.30'
$ws.Range('E51').Value = '  -2.44%  '
